$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Rename the sheets
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws1.Name = "row format"
$ws2.Name = "in column format"

# ------------------------------------------------------------------
# 2. Fill in the "in column format" sheet (Sheet2) with the db column
#    names, one per row, taken from the same text that already lives
#    on "row format" (Sheet1) row 1.
# ------------------------------------------------------------------
$values = @{
    2  = "magicPwd"
    3  = "clientStatus"
    4  = "toolStatus"
    5  = "TBD1"
    6  = "TBD2"
    7  = "TBD3"
    8  = "TBD4"
    9  = "bisiName"
    10 = "bisiStatus"
    11 = "bisiStartDate"
    12 = "bisiEndDate"
    13 = "bisiSumAssured"
    14 = "bisiComission"
    15 = "bisiPrevMonthComission"
    16 = "bisiComissionHistoryData"
    17 = "bisiPplList"
    18 = "bisiTotalPpl"
    19 = "bisiTotalMonths"
    20 = "TBD7"
    21 = "TBD8"
    22 = "TBD9"
    23 = "TBD10"
    24 = "personName"
    25 = "personDob"
    26 = "personAadhar"
    27 = "personPhone"
    28 = "personPhone"
    29 = "personAssociatedBisisList"
    30 = "personBisiEncashStatus"
    31 = "personBisiNameWhichIsEncashed"
    32 = "personBisiEncashedValue"
    33 = "TBD11"
    34 = "TBD12"
    35 = "TBD13"
    36 = "TBD14"
    37 = "TBD15"
}

foreach ($r in 2..37) {
    $ws2.Range("A$r").Value = $values[$r]
}

# "add bisi" marker column next to the bisi fields that are editable
foreach ($r in 9,10,11,12,13,18,19) {
    $ws2.Range("B$r").Value = "add bisi"
}

# ------------------------------------------------------------------
# 3. Formatting: highlight the "bisi" input rows in yellow, make the
#    two new totals rows bold+yellow, and the read-only bisi rows a
#    plain (theme) white fill.
# ------------------------------------------------------------------
foreach ($r in 9,10,11,12,13) {
    $ws2.Range("A$r").Interior.Color = 65535
}

foreach ($r in 18,19) {
    $ws2.Range("A$r").Interior.Color = 65535
    $ws2.Range("A$r").Font.Bold = $true
}

foreach ($r in 14,15,16,17) {
    $ws2.Range("A$r").Interior.Color = 65535
    $ws2.Range("A$r").Interior.ThemeColor = 2
}

# ------------------------------------------------------------------
# 4. Column width / page setup for the "in column format" sheet
# ------------------------------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 25.893229166666668
$ws2.PageSetup.Orientation = 1

# ------------------------------------------------------------------
# 5. View state: Sheet2 ("in column format") becomes the active /
#    selected tab, scrolled so row 4 is at the top, with A19 selected.
# ------------------------------------------------------------------
$ws1.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$ws1.Range("G15").Select()

$ws2.Activate()
$win2 = $excel.ActiveWindow
$win2.ScrollRow = 4
$ws2.Range("A19").Select()

Write-Output "done"
